$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.391.97"
$ws.Range("E2").Value = "  -3.91%  "
$ws.Range("D3").Value = "2.988.10"
$ws.Range("E3").Value = "  -5.39%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "565.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "125.80"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.92%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "2.982.29"
$ws.Range("E8").Value = "  -5.44%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.496"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.55%  "
$ws.Range("E10").Value = "  -5.71%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.07"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.01%  "
$ws.Range("E12").Value = "  -3.55%  "
$ws.Range("E13").Value = "  -5.83%  "
$ws.Range("E14").Value = "  -6.12%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.119"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.25%  "
$ws.Range("D16").Value = "3.480.61"
$ws.Range("E16").Value = "  -5.40%  "
$ws.Range("D17").Value = "60.642.10"
$ws.Range("E17").Value = "  -3.46%  "
$ws.Range("D18").Value = "2.984.99"
$ws.Range("E18").Value = "  -5.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -6.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "429.91"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.58%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.01"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.02%  "
$ws.Range("E22").Value = "  -5.71%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.10"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.54%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.88%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "78.43"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.90%  "
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.48"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.87%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.89"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.61%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.20"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.02"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -10.33%  "
$ws.Range("E33").Value = "  -9.75%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.25"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.58%  "
$ws.Range("E35").Value = "  -8.51%  "
$ws.Range("E36").Value = "  -5.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "49.18"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.84%  "
$ws.Range("E38").Value = "  -6.02%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0357"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.70%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.78"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "373.57"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.49%  "
$ws.Range("E42").Value = "  -4.62%  "
$ws.Range("D43").Value = "2.663.45"
$ws.Range("E43").Value = "  -4.77%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.41"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.22%  "
$ws.Range("E46").Value = "  -6.50%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "118.82"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.91%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.95"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -7.23%  "
$ws.Range("E49").Value = "  -4.82%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "32.48"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.90%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.29"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.89%  "
